$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report generation: the Overview / zh-cn / de-de sheets get their
# status updated from "Ready for handoff" to "Handed back: in sync with
# en-US", the zh-cn locale's target/handback info is filled in, and the
# de-de locale's target/handback info is filled in as well.
# ---------------------------------------------------------------------------

$statusText = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# Widen the zh-cn/de-de status columns on the Overview sheet to fit the
# longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.1667
$overview.Columns.Item(6).ColumnWidth = 29.1667

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("I2").Value = "a.md"
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-10-17 15:40:22"

$zhcn.Range("I3").Value = "a.md"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-10-17 15:40:22"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2f6d5ac46497618ff02a23baca4756261bd7666/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2f6d5ac46497618ff02a23baca4756261bd7666/e2e/a.md", "", "", "a.md")

$zhcn.Columns.Item(3).ColumnWidth = 29.1667
$zhcn.Columns.Item(10).ColumnWidth = 39.1667

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("I2").Value = "a.md"
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-10-17 15:41:00"

$dede.Range("I3").Value = "a.md"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-10-17 15:41:00"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2f6d5ac46497618ff02a23baca4756261bd7666/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2f6d5ac46497618ff02a23baca4756261bd7666/e2e/a.md", "", "", "a.md")

$dede.Columns.Item(3).ColumnWidth = 29.1667
$dede.Columns.Item(10).ColumnWidth = 39.1667
